# Fruta / hortaliza, semanal
#
# Weekly data refresh: a new price observation (row 9, fecha 2023-01-06)
# is inserted into the "Arándano (blue)" sheet, pushing the existing
# rows 9-11 down to 10-12 (their contents are unchanged, only their row
# position shifts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 9; everything from the old row 9
# downward (old 9 -> new 10, old 10 -> new 11, old 11 -> new 12) shifts
# down automatically and keeps its existing values/formatting.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with this week's observation.
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C9").Value = 'Ñuble'
$ws.Range("D9").Value = 44932
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 'Fruta'
$ws.Range("G9").Value = 100101
$ws.Range("H9").Value = 'Berries'
$ws.Range("I9").Value = 100101001
$ws.Range("J9").Value = 'Arándano (blue)'
$ws.Range("K9").Value = 'Sin especificar'
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 3000
$ws.Range("O9").Value = 3000
$ws.Range("P9").Value = 3000
$ws.Range("Q9").Value = '$/bandeja 2 kilos'
$ws.Range("R9").Value = 'Provincia de Diguillín'
$ws.Range("S9").Value = 1500
$ws.Range("T9").Value = 2
